# Sync automatico del tracker: marca como "Completed" las predicciones de los
# partidos ya finalizados (filas 249-281) y registra resultado real, profit,
# ROI y marca de tiempo del envio.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Predictions")

$ws.Range("L249").Value = "Completed"
$ws.Range("M249").Value = "Away Win"
$ws.Range("N249").Value = "Fallo"
$ws.Range("O249").Value = -1.7
$ws.Range("P249").Value = -100
$ws.Range("Q249").Value = "2025-10-06 04:26:46"

$ws.Range("L250").Value = "Completed"
$ws.Range("M250").Value = "Draw"
$ws.Range("N250").Value = "Fallo"
$ws.Range("O250").Value = -2.9
$ws.Range("P250").Value = -100
$ws.Range("Q250").Value = "2025-10-06 04:26:46"

$ws.Range("L251").Value = "Completed"
$ws.Range("M251").Value = "Home Win"
$ws.Range("N251").Value = "Acierto"
$ws.Range("O251").Value = 1.45
$ws.Range("P251").Value = 50
$ws.Range("Q251").Value = "2025-10-06 04:26:46"

$ws.Range("L252").Value = "Completed"
$ws.Range("M252").Value = "Home Win"
$ws.Range("N252").Value = "Acierto"
$ws.Range("O252").Value = 1.3
$ws.Range("P252").Value = 45
$ws.Range("Q252").Value = "2025-10-06 04:26:46"

$ws.Range("L253").Value = "Completed"
$ws.Range("M253").Value = "Home Win"
$ws.Range("N253").Value = "Acierto"
$ws.Range("O253").Value = 1.76
$ws.Range("P253").Value = 80
$ws.Range("Q253").Value = "2025-10-06 04:26:46"

$ws.Range("L254").Value = "Completed"
$ws.Range("M254").Value = "Home Win"
$ws.Range("N254").Value = "Acierto"
$ws.Range("O254").Value = 1.65
$ws.Range("P254").Value = 57
$ws.Range("Q254").Value = "2025-10-06 04:26:46"

$ws.Range("L255").Value = "Completed"
$ws.Range("M255").Value = "Draw"
$ws.Range("N255").Value = "Fallo"
$ws.Range("O255").Value = -0.9
$ws.Range("P255").Value = -100
$ws.Range("Q255").Value = "2025-10-06 04:26:46"

$ws.Range("L256").Value = "Completed"
$ws.Range("M256").Value = "Home Win"
$ws.Range("N256").Value = "Acierto"
$ws.Range("O256").Value = 1.44
$ws.Range("P256").Value = 60
$ws.Range("Q256").Value = "2025-10-06 04:26:46"

$ws.Range("L257").Value = "Completed"
$ws.Range("M257").Value = "Home Win"
$ws.Range("N257").Value = "Acierto"
$ws.Range("O257").Value = 1.39
$ws.Range("P257").Value = 48
$ws.Range("Q257").Value = "2025-10-06 04:26:46"

$ws.Range("L258").Value = "Completed"
$ws.Range("M258").Value = "Home Win"
$ws.Range("N258").Value = "Acierto"
$ws.Range("O258").Value = 1.56
$ws.Range("P258").Value = 65
$ws.Range("Q258").Value = "2025-10-06 04:26:46"

$ws.Range("L259").Value = "Completed"
$ws.Range("M259").Value = "Away Win"
$ws.Range("N259").Value = "Acierto"
$ws.Range("O259").Value = 0.55
$ws.Range("P259").Value = 138
$ws.Range("Q259").Value = "2025-10-06 04:26:46"

$ws.Range("L260").Value = "Completed"
$ws.Range("M260").Value = "Home Win"
$ws.Range("N260").Value = "Acierto"
$ws.Range("O260").Value = 1.62
$ws.Range("P260").Value = 60
$ws.Range("Q260").Value = "2025-10-06 04:26:46"

$ws.Range("L261").Value = "Completed"
$ws.Range("M261").Value = "Home Win"
$ws.Range("N261").Value = "Acierto"
$ws.Range("O261").Value = 1.59
$ws.Range("P261").Value = 55
$ws.Range("Q261").Value = "2025-10-06 04:26:46"

$ws.Range("L262").Value = "Completed"
$ws.Range("M262").Value = "Draw"
$ws.Range("N262").Value = "Fallo"
$ws.Range("O262").Value = -1.2
$ws.Range("P262").Value = -100
$ws.Range("Q262").Value = "2025-10-06 04:26:46"

$ws.Range("L263").Value = "Completed"
$ws.Range("M263").Value = "Away Win"
$ws.Range("N263").Value = "Fallo"
$ws.Range("O263").Value = -1.6
$ws.Range("P263").Value = -100
$ws.Range("Q263").Value = "2025-10-06 04:26:46"

$ws.Range("L264").Value = "Completed"
$ws.Range("M264").Value = "Home Win"
$ws.Range("N264").Value = "Acierto"
$ws.Range("O264").Value = 1.22
$ws.Range("P264").Value = 42
$ws.Range("Q264").Value = "2025-10-06 04:26:46"

$ws.Range("L265").Value = "Completed"
$ws.Range("M265").Value = "Home Win"
$ws.Range("N265").Value = "Fallo"
$ws.Range("O265").Value = -2.9
$ws.Range("P265").Value = -100
$ws.Range("Q265").Value = "2025-10-06 04:26:46"

$ws.Range("L266").Value = "Completed"
$ws.Range("M266").Value = "Draw"
$ws.Range("N266").Value = "Fallo"
$ws.Range("O266").Value = -1
$ws.Range("P266").Value = -100
$ws.Range("Q266").Value = "2025-10-06 04:26:46"

$ws.Range("L267").Value = "Completed"
$ws.Range("M267").Value = "Draw"
$ws.Range("N267").Value = "Fallo"
$ws.Range("O267").Value = -0.2
$ws.Range("P267").Value = -100
$ws.Range("Q267").Value = "2025-10-06 04:26:46"

$ws.Range("L268").Value = "Completed"
$ws.Range("M268").Value = "Draw"
$ws.Range("N268").Value = "Fallo"
$ws.Range("O268").Value = -2
$ws.Range("P268").Value = -100
$ws.Range("Q268").Value = "2025-10-06 04:26:46"

$ws.Range("L269").Value = "Completed"
$ws.Range("M269").Value = "Home Win"
$ws.Range("N269").Value = "Acierto"
$ws.Range("O269").Value = 1.61
$ws.Range("P269").Value = 70
$ws.Range("Q269").Value = "2025-10-06 15:23:53"

$ws.Range("L270").Value = "Completed"
$ws.Range("M270").Value = "Away Win"
$ws.Range("N270").Value = "Acierto"
$ws.Range("O270").Value = 1.69
$ws.Range("P270").Value = 65
$ws.Range("Q270").Value = "2025-10-06 15:23:53"

$ws.Range("L271").Value = "Completed"
$ws.Range("M271").Value = "Home Win"
$ws.Range("N271").Value = "Acierto"
$ws.Range("O271").Value = 1.28
$ws.Range("P271").Value = 44
$ws.Range("Q271").Value = "2025-10-06 15:23:53"

$ws.Range("L272").Value = "Completed"
$ws.Range("M272").Value = "Home Win"
$ws.Range("N272").Value = "Acierto"
$ws.Range("O272").Value = 0.6899999999999999
$ws.Range("P272").Value = 115
$ws.Range("Q272").Value = "2025-10-06 15:23:53"

$ws.Range("L273").Value = "Completed"
$ws.Range("M273").Value = "Away Win"
$ws.Range("N273").Value = "Fallo"
$ws.Range("O273").Value = -1.1
$ws.Range("P273").Value = -100
$ws.Range("Q273").Value = "2025-10-06 15:23:53"

$ws.Range("L274").Value = "Completed"
$ws.Range("M274").Value = "Home Win"
$ws.Range("N274").Value = "Acierto"
$ws.Range("O274").Value = 1.42
$ws.Range("P274").Value = 75
$ws.Range("Q274").Value = "2025-10-06 15:23:53"

$ws.Range("L275").Value = "Completed"
$ws.Range("M275").Value = "Draw"
$ws.Range("N275").Value = "Fallo"
$ws.Range("O275").Value = -2.4
$ws.Range("P275").Value = -100
$ws.Range("Q275").Value = "2025-10-06 15:23:53"

$ws.Range("L276").Value = "Completed"
$ws.Range("M276").Value = "Home Win"
$ws.Range("N276").Value = "Acierto"
$ws.Range("O276").Value = 1.45
$ws.Range("P276").Value = 85
$ws.Range("Q276").Value = "2025-10-06 15:23:53"

$ws.Range("L277").Value = "Completed"
$ws.Range("M277").Value = "Draw"
$ws.Range("N277").Value = "Fallo"
$ws.Range("O277").Value = -2.9
$ws.Range("P277").Value = -100
$ws.Range("Q277").Value = "2025-10-06 15:23:53"

$ws.Range("L278").Value = "Completed"
$ws.Range("M278").Value = "Draw"
$ws.Range("N278").Value = "Fallo"
$ws.Range("O278").Value = -2.3
$ws.Range("P278").Value = -100
$ws.Range("Q278").Value = "2025-10-06 15:23:53"

$ws.Range("L279").Value = "Completed"
$ws.Range("M279").Value = "Draw"
$ws.Range("N279").Value = "Fallo"
$ws.Range("O279").Value = -2.2
$ws.Range("P279").Value = -100
$ws.Range("Q279").Value = "2025-10-06 15:23:53"

$ws.Range("L280").Value = "Completed"
$ws.Range("M280").Value = "Draw"
$ws.Range("N280").Value = "Fallo"
$ws.Range("O280").Value = -1.4
$ws.Range("P280").Value = -100
$ws.Range("Q280").Value = "2025-10-06 15:23:53"

$ws.Range("L281").Value = "Completed"
$ws.Range("M281").Value = "Home Win"
$ws.Range("N281").Value = "Acierto"
$ws.Range("O281").Value = 1.68
$ws.Range("P281").Value = 60
$ws.Range("Q281").Value = "2025-10-06 15:23:53"
